# Updated cryptos list with GitHub Actions
# Applies per-cell text updates (price + volume% columns, plus a row
# 32/33 coin swap) as captured by the source diff. Numeric-looking
# "Price" strings are written through a text NumberFormat so Excel
# keeps them as literal text (matching the inline-string storage in
# the target file) instead of silently parsing them into numbers,
# then ClearFormats() restores the cell to its original (default)
# style so only the value itself changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.583.90'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -2.29%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.816.97'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.88%  '
$ws.Range("E4").Value = '  +0.99%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.009'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '309.06'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.46%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4580'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.48%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3672'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07164'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8795'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07804'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.45%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.44'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -3.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.789.33'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.302'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.45%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.392'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.01%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '86.34'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -4.88%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.011'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.99%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008609'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -3.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.008'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.642.50'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.32'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.74%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.019'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.48'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.987'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.17'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.38%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.03'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.081'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '113.23'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.29%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.871'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -3.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08695'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.058'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.57%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.495'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.20%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7351'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -4.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.121'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.005'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.573'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -5.88%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.082'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01941'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05120'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.49%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.906'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.03%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.009'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5043'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1564'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.88%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.185'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.010'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4642'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.27%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.03'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '101.06'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.77%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06039'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '64.34'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.85%  '
